$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    9  = -10.85
    13 = -13.311
    16 = -13.026
    18 = -11.62
    20 = -12.578
    26 = -12.494
    27 = -12.769
    29 = -12.041
    35 = -12.377
    36 = -12.428
    45 = -12.758
    55 = -13.381
    57 = -13.426
    69 = -10.537
    76 = -13.201
    78 = -12.5
    82 = -11.888
    83 = -13.085
    93 = -11.587
    97 = -12.028
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
